$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook-level settings ---
# Switch the reference style to R1C1 (the commit's calcPr refMode="R1C1" change)
$excel.ReferenceStyle = 2

# --- Worksheet header row changes ---
# Insert two new blank columns where "撥款日期" (R) and "到期日" (S) will go,
# pushing the existing "檢核訊息" column from R to T.
[void]$ws.Range("R1:S1").EntireColumn.Insert()

# The old " 加碼值" header (now in P1) is renamed to "合約加碼值".
$ws.Range("P1").Value = "合約加碼值"

# Fill in the two newly-inserted header cells.
$ws.Range("R1").Value = "撥款日期"
$ws.Range("S1").Value = "到期日"

# Give the new header cells the same center-aligned header formatting used by
# their neighboring header cells (matches the shared style already used by
# columns such as C1/D1/E1/I1/J1/K1) instead of the plain format Insert() guessed.
[void]$ws.Range("C1").Copy()
[void]$ws.Range("R1:S1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Match the new columns' width to the adjacent "擬調利率" column (Q).
$ws.Range("R1:S1").ColumnWidth = $ws.Range("Q1").ColumnWidth

# Reflect the saved selection state: the whole header row selected.
[void]$ws.Rows("1:1").Select()
